$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Populate cell values ---
# The fill order below is deliberate: it reproduces the shared-string table
# order of the target workbook (header row first - left to right, skipping
# the later-added "Sitename" column - then the original data columns in
# their original authoring order, and finally the newly mapped "Sitename"
# column E last - matching "map IRIS" in the commit message).

# Header row (columns A-D, F-K; E/"Sitename" header added later below)
$ws.Range("A1").Value = "StudyId"
$ws.Range("B1").Value = "StudySiteId"
$ws.Range("C1").Value = "IRBAgency"
$ws.Range("D1").Value = "IRBNumber"
$ws.Range("F1").Value = "EventID"
$ws.Range("G1").Value = "Event"
$ws.Range("H1").Value = "EventCreationDate"
$ws.Range("I1").Value = "EventOutcome"
$ws.Range("J1").Value = "TaskCompletionDate"
$ws.Range("K1").Value = "EventCompletionDate"

# Data rows 2-3, original columns
$ws.Range("C2").Value = "BRANY"
$ws.Range("C3").Value = "BRANY"
$ws.Range("F2").Value = "Einstein"
$ws.Range("F3").Value = "Einstein"
$ws.Range("A2").Value = """654654"""
$ws.Range("A3").Value = """4564654"""
$ws.Range("B2").Value = """632132"""
$ws.Range("B3").Value = """64654"""
$ws.Range("D2").Value = """12345"""
$ws.Range("D3").Value = """123465498"""
$ws.Range("G2").Value = """IRBAmendment"""
$ws.Range("G3").Value = """IRBAmendment"""
$ws.Range("H2").Value = "2015-03-12T08:22:12"
$ws.Range("H3").Value = "2015-03-12T08:22:12"
$ws.Range("J3").Value = "2015-04-12T08:22:12"
$ws.Range("K2").Value = "2015-04-12T08:22:12"

# Newly mapped IRIS site-name column (E)
$ws.Range("E1").Value = "Sitename"
$ws.Range("E2").Value = "Einstein-Montefiore"
$ws.Range("E3").Value = "Einstein-Montefiore"

# --- Column widths (best-fit sizing after the new data/column were added) ---
$ws.Range("B1").ColumnWidth = 13.833333333333334
$ws.Range("C1").ColumnWidth = 11.333333333333334
$ws.Range("D1").ColumnWidth = 12.0
$ws.Range("E1").ColumnWidth = 12.0
$ws.Range("G1").ColumnWidth = 10.0
$ws.Range("H1").ColumnWidth = 19.333333333333332
$ws.Range("I1").ColumnWidth = 19.833333333333332
$ws.Range("J1").ColumnWidth = 20.5
$ws.Range("K1").ColumnWidth = 20.666666666666668

# --- Selection ends up on E9 (as in the saved file) ---
$ws.Range("E9").Select() | Out-Null
